$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.293.11"
$ws.Range("E2").Value = "'  +1.73%  "

$ws.Range("D3").Value = "'2.653.56"
$ws.Range("E3").Value = "'  +2.26%  "

$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'580.38"
$ws.Range("E5").Value = "'  +0.42%  "

$ws.Range("D6").Value = "'144.54"
$ws.Range("E6").Value = "'  +1.21%  "

$ws.Range("E7").Value = "'  -0.07%  "

$ws.Range("E8").Value = "'  +0.41%  "

$ws.Range("D9").Value = "'6.58"
$ws.Range("E9").Value = "'  +0.87%  "

$ws.Range("E10").Value = "'  +4.46%  "

$ws.Range("D11").Value = "'0.384"
$ws.Range("E11").Value = "'  +3.38%  "

$ws.Range("E12").Value = "'  +0.88%  "

$ws.Range("D13").Value = "'3.121.98"
$ws.Range("E13").Value = "'  +2.22%  "

$ws.Range("D14").Value = "'26.24"
$ws.Range("E14").Value = "'  +7.17%  "

$ws.Range("D15").Value = "'61.206.43"
$ws.Range("E15").Value = "'  +1.58%  "

$ws.Range("D16").Value = "'0.0000146"
$ws.Range("E16").Value = "'  +3.65%  "

$ws.Range("D17").Value = "'2.660.86"
$ws.Range("E17").Value = "'  +2.31%  "

$ws.Range("D18").Value = "'11.68"
$ws.Range("E18").Value = "'  +2.01%  "

$ws.Range("D19").Value = "'4.79"

$ws.Range("D20").Value = "'355.14"
$ws.Range("E20").Value = "'  +2.38%  "

$ws.Range("D21").Value = "'6.88"
$ws.Range("E21").Value = "'  -0.29%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "'  -0.24%  "

$ws.Range("E23").Value = "'  +1.32%  "

$ws.Range("D24").Value = "'64.61"
$ws.Range("E24").Value = "'  +2.53%  "

$ws.Range("E25").Value = "'  +3.05%  "

$ws.Range("D26").Value = "'8.51"
$ws.Range("E26").Value = "'  +6.21%  "

$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "'  -0.16%  "

$ws.Range("E28").Value = "'  +8.17%  "

$ws.Range("D29").Value = "'0.0₃0817"
$ws.Range("E29").Value = "'  +3.34%  "

$ws.Range("D30").Value = "'6.89"
$ws.Range("E30").Value = "'  +8.06%  "

$ws.Range("D31").Value = "'168.87"
$ws.Range("E31").Value = "'  +2.46%  "

$ws.Range("E32").Value = "'  -0.05%  "

$ws.Range("D33").Value = "'20.19"
$ws.Range("E33").Value = "'  +3.95%  "

$ws.Range("D34").Value = "'1.13"
$ws.Range("E34").Value = "'  +14.88%  "

$ws.Range("D35").Value = "'4.69"
$ws.Range("E35").Value = "'  +9.41%  "

$ws.Range("E36").Value = "'  +10.34%  "

$ws.Range("D37").Value = "'0.984"
$ws.Range("E37").Value = "'  +17.21%  "

$ws.Range("E38").Value = "'  +6.14%  "

$ws.Range("D39").Value = "'337.77"
$ws.Range("E39").Value = "'  +8.91%  "

$ws.Range("D40").Value = "'4.15"
$ws.Range("E40").Value = "'  +6.30%  "

$ws.Range("D41").Value = "'38.49"
$ws.Range("E41").Value = "'  +1.09%  "

$ws.Range("E42").Value = "'  +6.27%  "

$ws.Range("D43").Value = "'0.0580"
$ws.Range("E43").Value = "'  +5.84%  "

$ws.Range("D44").Value = "'20.67"
$ws.Range("E44").Value = "'  +4.69%  "

$ws.Range("B45").Value = "'InjectiveProtocol"
$ws.Range("C45").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'21.12"
$ws.Range("E45").Value = "'  +5.48%  "

$ws.Range("B46").Value = "'VeChain"
$ws.Range("C46").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0254"
$ws.Range("E46").Value = "'  +5.54%  "

$ws.Range("D47").Value = "'135.77"
$ws.Range("E47").Value = "'  +0.55%  "

$ws.Range("B48").Value = "'Mantle"
$ws.Range("C48").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.627"
$ws.Range("E48").Value = "'  +4.07%  "

$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "'  +1.05%  "

$ws.Range("D50").Value = "'0.996"
$ws.Range("E50").Value = "'  -0.25%  "

$ws.Range("D51").Value = "'2.095.60"
$ws.Range("E51").Value = "'  +3.79%  "
